$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, as scraped from the source diff.
# D-column values are price strings that must stay as literal text
# (they may contain multiple "." thousand separators or significant
# trailing/leading zeros that a numeric auto-conversion would destroy),
# so they are written with a leading apostrophe to force text, and the
# cell's original Style is restored afterwards so no formatting changes
# are introduced.
$updates = @(
    @{ Ref = 'D2'; Value = '27.962.98'; ForceText = $true }
    @{ Ref = 'E2'; Value = '  +2.27%  '; ForceText = $false }
    @{ Ref = 'D3'; Value = '1.813.09'; ForceText = $true }
    @{ Ref = 'E3'; Value = '  +1.56%  '; ForceText = $false }
    @{ Ref = 'D4'; Value = '0.9980'; ForceText = $true }
    @{ Ref = 'E4'; Value = '  -0.71%  '; ForceText = $false }
    @{ Ref = 'D5'; Value = '337.48'; ForceText = $true }
    @{ Ref = 'D6'; Value = '0.9964'; ForceText = $true }
    @{ Ref = 'E6'; Value = '  -0.72%  '; ForceText = $false }
    @{ Ref = 'D7'; Value = '0.3927'; ForceText = $true }
    @{ Ref = 'E7'; Value = '  +3.90%  '; ForceText = $false }
    @{ Ref = 'D8'; Value = '0.3491'; ForceText = $true }
    @{ Ref = 'E8'; Value = '  +2.13%  '; ForceText = $false }
    @{ Ref = 'D9'; Value = '48.30'; ForceText = $true }
    @{ Ref = 'E9'; Value = '  +0.52%  '; ForceText = $false }
    @{ Ref = 'D10'; Value = '1.203'; ForceText = $true }
    @{ Ref = 'E10'; Value = '  +0.39%  '; ForceText = $false }
    @{ Ref = 'D11'; Value = '0.07585'; ForceText = $true }
    @{ Ref = 'E11'; Value = '  +1.95%  '; ForceText = $false }
    @{ Ref = 'D12'; Value = '0.9957'; ForceText = $true }
    @{ Ref = 'E12'; Value = '  -0.69%  '; ForceText = $false }
    @{ Ref = 'D13'; Value = '22.19'; ForceText = $true }
    @{ Ref = 'E13'; Value = '  +1.39%  '; ForceText = $false }
    @{ Ref = 'D14'; Value = '6.530'; ForceText = $true }
    @{ Ref = 'E14'; Value = '  +1.22%  '; ForceText = $false }
    @{ Ref = 'D15'; Value = '1.812.92'; ForceText = $true }
    @{ Ref = 'E15'; Value = '  +1.52%  '; ForceText = $false }
    @{ Ref = 'D16'; Value = '7.206'; ForceText = $true }
    @{ Ref = 'E16'; Value = '  +2.73%  '; ForceText = $false }
    @{ Ref = 'D17'; Value = '0.00001108'; ForceText = $true }
    @{ Ref = 'D18'; Value = '0.06678'; ForceText = $true }
    @{ Ref = 'D19'; Value = '85.24'; ForceText = $true }
    @{ Ref = 'E19'; Value = '  +1.12%  '; ForceText = $false }
    @{ Ref = 'D20'; Value = '0.9965'; ForceText = $true }
    @{ Ref = 'E20'; Value = '  -0.74%  '; ForceText = $false }
    @{ Ref = 'D21'; Value = '17.90'; ForceText = $true }
    @{ Ref = 'E21'; Value = '  +3.43%  '; ForceText = $false }
    @{ Ref = 'D22'; Value = '6.580'; ForceText = $true }
    @{ Ref = 'E22'; Value = '  +2.20%  '; ForceText = $false }
    @{ Ref = 'D23'; Value = '27.951.56'; ForceText = $true }
    @{ Ref = 'E23'; Value = '  +2.40%  '; ForceText = $false }
    @{ Ref = 'D24'; Value = '12.88'; ForceText = $true }
    @{ Ref = 'E24'; Value = '  +2.96%  '; ForceText = $false }
    @{ Ref = 'D25'; Value = '2.407'; ForceText = $true }
    @{ Ref = 'E25'; Value = '  -1.84%  '; ForceText = $false }
    @{ Ref = 'D26'; Value = '2.561'; ForceText = $true }
    @{ Ref = 'E26'; Value = '  +0.60%  '; ForceText = $false }
    @{ Ref = 'D27'; Value = '1.480'; ForceText = $true }
    @{ Ref = 'E27'; Value = '  +0.66%  '; ForceText = $false }
    @{ Ref = 'D28'; Value = '21.36'; ForceText = $true }
    @{ Ref = 'E28'; Value = '  +0.45%  '; ForceText = $false }
    @{ Ref = 'D29'; Value = '155.06'; ForceText = $true }
    @{ Ref = 'E29'; Value = '  +3.38%  '; ForceText = $false }
    @{ Ref = 'D30'; Value = '2.017.66'; ForceText = $true }
    @{ Ref = 'E30'; Value = '  +1.51%  '; ForceText = $false }
    @{ Ref = 'D31'; Value = '136.09'; ForceText = $true }
    @{ Ref = 'E31'; Value = '  +2.43%  '; ForceText = $false }
    @{ Ref = 'D32'; Value = '4.036'; ForceText = $true }
    @{ Ref = 'E32'; Value = '  -0.60%  '; ForceText = $false }
    @{ Ref = 'D33'; Value = '6.161'; ForceText = $true }
    @{ Ref = 'E33'; Value = '  +1.16%  '; ForceText = $false }
    @{ Ref = 'D34'; Value = '0.08853'; ForceText = $true }
    @{ Ref = 'E34'; Value = '  +2.83%  '; ForceText = $false }
    @{ Ref = 'D35'; Value = '13.33'; ForceText = $true }
    @{ Ref = 'E35'; Value = '  +1.23%  '; ForceText = $false }
    @{ Ref = 'D36'; Value = '5.559'; ForceText = $true }
    @{ Ref = 'E36'; Value = '  +2.77%  '; ForceText = $false }
    @{ Ref = 'D37'; Value = '0.02432'; ForceText = $true }
    @{ Ref = 'E37'; Value = '  +3.95%  '; ForceText = $false }
    @{ Ref = 'D38'; Value = '0.6936'; ForceText = $true }
    @{ Ref = 'E38'; Value = '  +1.15%  '; ForceText = $false }
    @{ Ref = 'D39'; Value = '0.06551'; ForceText = $true }
    @{ Ref = 'E39'; Value = '  +3.27%  '; ForceText = $false }
    @{ Ref = 'E40'; Value = '  -3.23%  '; ForceText = $false }
    @{ Ref = 'D41'; Value = '0.2230'; ForceText = $true }
    @{ Ref = 'E41'; Value = '  +1.91%  '; ForceText = $false }
    @{ Ref = 'D42'; Value = '1.268'; ForceText = $true }
    @{ Ref = 'E42'; Value = '  -0.25%  '; ForceText = $false }
    @{ Ref = 'D43'; Value = '8.581'; ForceText = $true }
    @{ Ref = 'E43'; Value = '  -2.33%  '; ForceText = $false }
    @{ Ref = 'D44'; Value = '14.70'; ForceText = $true }
    @{ Ref = 'E44'; Value = '  +2.35%  '; ForceText = $false }
    @{ Ref = 'D45'; Value = '0.6559'; ForceText = $true }
    @{ Ref = 'E45'; Value = '  +2.22%  '; ForceText = $false }
    @{ Ref = 'D46'; Value = '0.9960'; ForceText = $true }
    @{ Ref = 'E46'; Value = '  -0.67%  '; ForceText = $false }
    @{ Ref = 'D47'; Value = '3.865'; ForceText = $true }
    @{ Ref = 'E47'; Value = '  +0.53%  '; ForceText = $false }
    @{ Ref = 'D48'; Value = '2.167'; ForceText = $true }
    @{ Ref = 'E48'; Value = '  +2.87%  '; ForceText = $false }
    @{ Ref = 'D49'; Value = '132.47'; ForceText = $true }
    @{ Ref = 'E49'; Value = '  +2.58%  '; ForceText = $false }
    @{ Ref = 'D50'; Value = '0.07208'; ForceText = $true }
    @{ Ref = 'E50'; Value = '  +0.47%  '; ForceText = $false }
    @{ Ref = 'D51'; Value = '80.83'; ForceText = $true }
    @{ Ref = 'E51'; Value = '  +2.42%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.ForceText) {
        $origStyle = $cell.Style
        $cell.Value = "'" + $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}
